# Update the NATMI Ctf1-Lifr LR-pair sheet with recomputed TPM values.
#
# The new run drops "ECs" as a *sending* cluster (it remains valid as a
# target cluster) and keeps only "FAPs" and "MuSCs" as senders, each paired
# with all three target clusters (ECs, FAPs, MuSCs). That is 2 x 3 = 6 data
# rows (rows 2-7) instead of the previous 3 x 3 = 9 (rows 2-10), so the
# trailing three rows (8-10) are removed and the sheet dimension shrinks
# from A1:T10 to A1:T7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @("FAPs",  "Ctf1", "Lifr", "ECs",   2, 0.6666666666666666, 0.7996153333333332, 2.398846,           0.373961854065697,   0.3739618540656971, 3, 1, 34.53319033333333, 103.599571, 0.2461870921144496, 0.2461870921144496, 27.61326849945177,  248.519416495066, 0.09206458141416213, 0.09206458141416213),
    @("FAPs",  "Ctf1", "Lifr", "FAPs",  2, 0.6666666666666666, 0.7996153333333332, 2.398846,           0.373961854065697,   0.3739618540656971, 3, 1, 80.77474466666666, 242.324234, 0.575843103803214,  0.575843103803214,  64.58872438155154,  581.298519433964, 0.2153433547491956,  0.2153433547491956),
    @("FAPs",  "Ctf1", "Lifr", "MuSCs", 2, 0.6666666666666666, 0.7996153333333332, 2.398846,           0.373961854065697,   0.3739618540656971, 3, 1, 24.96420533333334, 74.892616,  0.1779698040823365, 0.1779698040823364, 19.96176136901511,  179.655852321136, 0.0665539179023394,  0.0665539179023394),
    @("MuSCs", "Ctf1", "Lifr", "ECs",   3, 1,                  1.338611666666667, 4.015835,           0.6260381459343028,  0.626038145934303,  3, 1, 34.53319033333333, 103.599571, 0.2461870921144496, 0.2461870921144496, 46.22653146742056,  416.038783206785, 0.1541225107002875,  0.1541225107002875),
    @("MuSCs", "Ctf1", "Lifr", "FAPs",  3, 1,                  1.338611666666667, 4.015835,           0.6260381459343028,  0.626038145934303,  3, 1, 80.77474466666666, 242.324234, 0.575843103803214,  0.575843103803214,  108.1260155828211,  973.13414024539,  0.3604997490540184,  0.3604997490540184),
    @("MuSCs", "Ctf1", "Lifr", "MuSCs", 3, 1,                  1.338611666666667, 4.015835,           0.6260381459343028,  0.626038145934303,  3, 1, 24.96420533333334, 74.892616,  0.1779698040823365, 0.1779698040823364, 33.41737650826223,  300.75638857436,  0.111415886179997,   0.111415886179997)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = 2 + $i
    $rowVals = $newData[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $r).Value = $rowVals[$j]
    }
}

# Remove the now-unused trailing rows (old rows 8-10), shifting cells up and
# shrinking the sheet's used range from A1:T10 to A1:T7.
$ws.Range("A8:T10").Delete()
